# Weekly update: a new "Ajo" (garlic) price record for
# "Feria Lagunitas de Puerto Montt" is inserted as row 226, pushing the
# existing historical rows (226..306) down by one (to 227..307).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 226, shifting rows 226:306 down to 227:307.
$ws.Rows(226).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A226").Value = 4
$ws.Range("B226").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C226").Value = 'Los Lagos'
$ws.Range("D226").Value = 44795
$ws.Range("E226").Value = 10
$ws.Range("F226").Value = 100112003
$ws.Range("G226").Value = 'Ajo'
$ws.Range("H226").Value = 'Chino'
$ws.Range("I226").Value = 'Primera'
$ws.Range("J226").Value = 70
$ws.Range("K226").Value = 26000
$ws.Range("L226").Value = 26000
$ws.Range("M226").Value = 26000
$ws.Range("N226").Value = '$/caja 10 kilos'
$ws.Range("O226").Value = 'China'
$ws.Range("P226").Value = 2600
$ws.Range("Q226").Value = 10
$ws.Range("R226").Value = 'Hortaliza'
